# Update "want to go" counts (column F) on the 展览 and 全部类型 sheets.
$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value  = 5
$ws1.Range("F3").Value  = 12738
$ws1.Range("F4").Value  = 25
$ws1.Range("F8").Value  = 15
$ws1.Range("F9").Value  = 5
$ws1.Range("F10").Value = 12635
$ws1.Range("F11").Value = 258
$ws1.Range("F12").Value = 13
$ws1.Range("F13").Value = 4940
$ws1.Range("F14").Value = 5891
$ws1.Range("F15").Value = 175
$ws1.Range("F20").Value = 5
$ws1.Range("F24").Value = 3
$ws1.Range("F25").Value = 81

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F3").Value  = 5
$ws4.Range("F4").Value  = 12738
$ws4.Range("F5").Value  = 25
$ws4.Range("F9").Value  = 15
$ws4.Range("F10").Value = 5
$ws4.Range("F11").Value = 12635
$ws4.Range("F12").Value = 258
$ws4.Range("F13").Value = 13
$ws4.Range("F14").Value = 4940
$ws4.Range("F15").Value = 5891
$ws4.Range("F16").Value = 175
$ws4.Range("F21").Value = 5
$ws4.Range("F25").Value = 3
$ws4.Range("F26").Value = 81

$wb.Save()
